$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(33, 8).Value = 349.37036
$ws.Cells.Item(33, 9).Value = 369.6087
$ws.Cells.Item(33, 10).Value = 233
$ws.Cells.Item(33, 11).Value = 369.6087
$ws.Cells.Item(33, 12).Value = 233
$ws.Cells.Item(33, 13).Value = -140.6087
$ws.Cells.Item(33, 14).Value = -691
$ws.Cells.Item(64, 8).Value = 3810.1428
$ws.Cells.Item(64, 9).Value = 3671.6667
$ws.Cells.Item(64, 10).Value = 4059.4
$ws.Cells.Item(64, 11).Value = 3671.6667
$ws.Cells.Item(64, 12).Value = 4059.4
$ws.Cells.Item(64, 13).Value = -3423.6667
$ws.Cells.Item(64, 14).Value = -4555.4
$ws.Cells.Item(67, 8).Value = 3810.1428
$ws.Cells.Item(67, 9).Value = 3671.6667
$ws.Cells.Item(67, 10).Value = 4059.4
$ws.Cells.Item(67, 11).Value = 3671.6667
$ws.Cells.Item(67, 12).Value = 4059.4
$ws.Cells.Item(67, 13).Value = -2813.6667
$ws.Cells.Item(67, 14).Value = -5775.4
$ws.Cells.Item(76, 8).Value = 3433.111
$ws.Cells.Item(76, 10).Value = 3819.6
$ws.Cells.Item(76, 12).Value = 3819.6
$ws.Cells.Item(76, 14).Value = -4449.6
$ws.Cells.Item(79, 8).Value = 3433.111
$ws.Cells.Item(79, 10).Value = 3819.6
$ws.Cells.Item(79, 12).Value = 3819.6
$ws.Cells.Item(79, 14).Value = -6003.6
$ws.Cells.Item(100, 8).Value = 1485.4445
$ws.Cells.Item(100, 9).Value = 1006.3571
$ws.Cells.Item(100, 10).Value = 2001.3846
$ws.Cells.Item(100, 11).Value = 1006.3571
$ws.Cells.Item(100, 12).Value = 2001.3846
$ws.Cells.Item(100, 13).Value = -465.3570999999999
$ws.Cells.Item(100, 14).Value = -3083.3846
$ws.Cells.Item(107, 8).Value = 10158.3
$ws.Cells.Item(107, 9).Value = 14386.857
$ws.Cells.Item(107, 10).Value = 291.66666
$ws.Cells.Item(107, 11).Value = 14386.857
$ws.Cells.Item(107, 12).Value = 291.66666
$ws.Cells.Item(107, 13).Value = -12466.857
$ws.Cells.Item(107, 14).Value = -4131.66666
$ws.Cells.Item(112, 8).Value = 1739.6207
$ws.Cells.Item(112, 10).Value = 1939.5416
$ws.Cells.Item(112, 12).Value = 5818.6248
$ws.Cells.Item(112, 14).Value = -8034.6248
$ws.Cells.Item(113, 8).Value = 3523.4783
$ws.Cells.Item(113, 9).Value = 2911.6667
$ws.Cells.Item(113, 10).Value = 4190.909
$ws.Cells.Item(113, 11).Value = 2911.6667
$ws.Cells.Item(113, 12).Value = 4190.909
$ws.Cells.Item(113, 13).Value = 342.3332999999998
$ws.Cells.Item(113, 14).Value = -10698.909
$ws.Cells.Item(116, 8).Value = 2559.9333
$ws.Cells.Item(116, 9).Value = 2318.0908
$ws.Cells.Item(116, 10).Value = 3225
$ws.Cells.Item(116, 11).Value = 2318.0908
$ws.Cells.Item(116, 12).Value = 3225
$ws.Cells.Item(116, 13).Value = 1123.9092
$ws.Cells.Item(116, 14).Value = -10109

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(5, 8).Value = 347.14285
$ws.Cells.Item(5, 9).Value = 334.44446
$ws.Cells.Item(5, 10).Value = 370
$ws.Cells.Item(5, 11).Value = 334.44446
$ws.Cells.Item(5, 12).Value = 370
$ws.Cells.Item(5, 13).Value = -222.44446
$ws.Cells.Item(5, 14).Value = -594
$ws.Cells.Item(45, 8).Value = 1421785.1
$ws.Cells.Item(45, 9).Value = 1819499
$ws.Cells.Item(45, 10).Value = 1378.5714
$ws.Cells.Item(45, 11).Value = 1819499
$ws.Cells.Item(45, 12).Value = 1378.5714
$ws.Cells.Item(45, 13).Value = -1819122
$ws.Cells.Item(45, 14).Value = -2132.5714
$ws.Cells.Item(97, 8).Value = 6529.7646
$ws.Cells.Item(97, 9).Value = 7427.4287
$ws.Cells.Item(97, 10).Value = 2340.6667
$ws.Cells.Item(97, 11).Value = 7427.4287
$ws.Cells.Item(97, 12).Value = 2340.6667
$ws.Cells.Item(97, 13).Value = -6931.4287
$ws.Cells.Item(97, 14).Value = -3332.6667

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(4, 8).Value = 347.14285
$ws.Cells.Item(4, 9).Value = 334.44446
$ws.Cells.Item(4, 10).Value = 370
$ws.Cells.Item(4, 11).Value = 334.44446
$ws.Cells.Item(4, 12).Value = 370
$ws.Cells.Item(4, 13).Value = -219.44446
$ws.Cells.Item(4, 14).Value = -600
$ws.Cells.Item(51, 8).Value = 47668
$ws.Cells.Item(51, 10).Value = 47668
$ws.Cells.Item(51, 12).Value = 47668
$ws.Cells.Item(51, 14).Value = -48650
$ws.Cells.Item(105, 8).Value = 3832.725
$ws.Cells.Item(105, 9).Value = 2028
$ws.Cells.Item(105, 11).Value = 2028
$ws.Cells.Item(105, 13).Value = -281
$ws.Cells.Item(134, 8).Value = 3532.6948
$ws.Cells.Item(134, 9).Value = 2199.2856
$ws.Cells.Item(134, 11).Value = 6597.8568
$ws.Cells.Item(134, 13).Value = -4062.8568

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(7, 8).Value = 153.8
$ws.Cells.Item(7, 9).Value = 58.916668
$ws.Cells.Item(7, 11).Value = 58.916668
$ws.Cells.Item(7, 13).Value = 54.083332
$ws.Cells.Item(31, 8).Value = 3638.7637
$ws.Cells.Item(31, 9).Value = 2170.4
$ws.Cells.Item(31, 10).Value = 5400.8
$ws.Cells.Item(31, 11).Value = 2170.4
$ws.Cells.Item(31, 12).Value = 5400.8
$ws.Cells.Item(31, 13).Value = -1875.4
$ws.Cells.Item(31, 14).Value = -5990.8
$ws.Cells.Item(34, 8).Value = 3638.7637
$ws.Cells.Item(34, 9).Value = 2170.4
$ws.Cells.Item(34, 10).Value = 5400.8
$ws.Cells.Item(34, 11).Value = 2170.4
$ws.Cells.Item(34, 12).Value = 5400.8
$ws.Cells.Item(34, 13).Value = -1968.4
$ws.Cells.Item(34, 14).Value = -5804.8
$ws.Cells.Item(62, 8).Value = 2598.5833
$ws.Cells.Item(62, 9).Value = 2369.2856
$ws.Cells.Item(62, 10).Value = 2919.6
$ws.Cells.Item(62, 11).Value = 2369.2856
$ws.Cells.Item(62, 12).Value = 2919.6
$ws.Cells.Item(62, 13).Value = -1745.2856
$ws.Cells.Item(62, 14).Value = -4167.6
$ws.Cells.Item(65, 8).Value = 2598.5833
$ws.Cells.Item(65, 9).Value = 2369.2856
$ws.Cells.Item(65, 10).Value = 2919.6
$ws.Cells.Item(65, 11).Value = 11846.428
$ws.Cells.Item(65, 12).Value = 14598
$ws.Cells.Item(65, 13).Value = -8726.428
$ws.Cells.Item(65, 14).Value = -20838
$ws.Cells.Item(80, 8).Value = 24416.6
$ws.Cells.Item(80, 10).Value = 24416.6
$ws.Cells.Item(80, 12).Value = 24416.6
$ws.Cells.Item(80, 14).Value = -26662.6
$ws.Cells.Item(83, 8).Value = 24416.6
$ws.Cells.Item(83, 10).Value = 24416.6
$ws.Cells.Item(83, 12).Value = 73249.79999999999
$ws.Cells.Item(83, 14).Value = -84481.79999999999

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(2, 8).Value = 868234.25
$ws.Cells.Item(2, 9).Value = 128.88889
$ws.Cells.Item(2, 10).Value = 1984369.8
$ws.Cells.Item(2, 11).Value = 773.33334
$ws.Cells.Item(2, 12).Value = 11906218.8
$ws.Cells.Item(2, 13).Value = -660.33334
$ws.Cells.Item(2, 14).Value = -11906444.8
$ws.Cells.Item(113, 8).Value = 29413066
$ws.Cells.Item(113, 9).Value = 66667884
$ws.Cells.Item(113, 10).Value = 1368.421
$ws.Cells.Item(113, 11).Value = 200003652
$ws.Cells.Item(113, 12).Value = 4105.263
$ws.Cells.Item(113, 13).Value = -200001482
$ws.Cells.Item(113, 14).Value = -8445.262999999999
$ws.Cells.Item(132, 8).Value = 1490.5
$ws.Cells.Item(132, 9).Value = 870.6
$ws.Cells.Item(132, 10).Value = 1933.2858
$ws.Cells.Item(132, 11).Value = 7835.400000000001
$ws.Cells.Item(132, 12).Value = 17399.5722
$ws.Cells.Item(132, 13).Value = -5305.400000000001
$ws.Cells.Item(132, 14).Value = -22459.5722
$ws.Cells.Item(134, 8).Value = 6564.4443
$ws.Cells.Item(134, 9).Value = 3616
$ws.Cells.Item(134, 10).Value = 10250
$ws.Cells.Item(134, 11).Value = 10848
$ws.Cells.Item(134, 12).Value = 30750
$ws.Cells.Item(134, 13).Value = -5778
$ws.Cells.Item(134, 14).Value = -40890
$ws.Cells.Item(138, 8).Value = 5758.2
$ws.Cells.Item(138, 9).Value = 963.5
$ws.Cells.Item(138, 10).Value = 12950.25
$ws.Cells.Item(138, 11).Value = 2890.5
$ws.Cells.Item(138, 12).Value = 38850.75
$ws.Cells.Item(138, 13).Value = 2249.5
$ws.Cells.Item(138, 14).Value = -49130.75

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(15, 8).Value = 19630
$ws.Cells.Item(15, 10).Value = 19630
$ws.Cells.Item(15, 12).Value = 19630
$ws.Cells.Item(15, 14).Value = -20206
$ws.Cells.Item(70, 8).Value = 12160.083
$ws.Cells.Item(70, 9).Value = 14308.632
$ws.Cells.Item(70, 10).Value = 3995.6
$ws.Cells.Item(70, 11).Value = 14308.632
$ws.Cells.Item(70, 12).Value = 3995.6
$ws.Cells.Item(70, 13).Value = -14038.632
$ws.Cells.Item(70, 14).Value = -4535.6
$ws.Cells.Item(73, 8).Value = 12160.083
$ws.Cells.Item(73, 9).Value = 14308.632
$ws.Cells.Item(73, 10).Value = 3995.6
$ws.Cells.Item(73, 11).Value = 14308.632
$ws.Cells.Item(73, 12).Value = 3995.6
$ws.Cells.Item(73, 13).Value = -13372.632
$ws.Cells.Item(73, 14).Value = -5867.6
$ws.Cells.Item(81, 8).Value = 19630
$ws.Cells.Item(81, 10).Value = 19630
$ws.Cells.Item(81, 12).Value = 19630
$ws.Cells.Item(81, 14).Value = -21626
$ws.Cells.Item(84, 8).Value = 19630
$ws.Cells.Item(84, 10).Value = 19630
$ws.Cells.Item(84, 12).Value = 58890
$ws.Cells.Item(84, 14).Value = -68874
$ws.Cells.Item(102, 8).Value = 3711.84
$ws.Cells.Item(102, 9).Value = 4621.7334
$ws.Cells.Item(102, 11).Value = 4621.7334
$ws.Cells.Item(102, 13).Value = -2999.7334

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(7, 8).Value = 9879.368
$ws.Cells.Item(7, 9).Value = 8601.429
$ws.Cells.Item(7, 10).Value = 10624.833
$ws.Cells.Item(7, 11).Value = 8601.429
$ws.Cells.Item(7, 12).Value = 10624.833
$ws.Cells.Item(7, 13).Value = -8489.429
$ws.Cells.Item(7, 14).Value = -10848.833
$ws.Cells.Item(126, 8).Value = 9879.368
$ws.Cells.Item(126, 9).Value = 8601.429
$ws.Cells.Item(126, 10).Value = 10624.833
$ws.Cells.Item(126, 11).Value = 25804.287
$ws.Cells.Item(126, 12).Value = 31874.499
$ws.Cells.Item(126, 13).Value = -23334.287
$ws.Cells.Item(126, 14).Value = -36814.499

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(62, 8).Value = 5371.875
$ws.Cells.Item(62, 9).Value = 5650
$ws.Cells.Item(62, 10).Value = 5205
$ws.Cells.Item(62, 11).Value = 5650
$ws.Cells.Item(62, 12).Value = 5205
$ws.Cells.Item(62, 13).Value = -5026
$ws.Cells.Item(62, 14).Value = -6453
$ws.Cells.Item(65, 8).Value = 5371.875
$ws.Cells.Item(65, 9).Value = 5650
$ws.Cells.Item(65, 10).Value = 5205
$ws.Cells.Item(65, 11).Value = 28250
$ws.Cells.Item(65, 12).Value = 26025
$ws.Cells.Item(65, 13).Value = -25130
$ws.Cells.Item(65, 14).Value = -32265
$ws.Cells.Item(75, 8).Value = 24900
$ws.Cells.Item(75, 9).Value = 10000
$ws.Cells.Item(75, 10).Value = 39800
$ws.Cells.Item(75, 11).Value = 10000
$ws.Cells.Item(75, 12).Value = 39800
$ws.Cells.Item(75, 13).Value = -9064
$ws.Cells.Item(75, 14).Value = -41672
$ws.Cells.Item(78, 8).Value = 24900
$ws.Cells.Item(78, 9).Value = 10000
$ws.Cells.Item(78, 10).Value = 39800
$ws.Cells.Item(78, 11).Value = 30000
$ws.Cells.Item(78, 12).Value = 119400
$ws.Cells.Item(78, 13).Value = -25320
$ws.Cells.Item(78, 14).Value = -128760
$ws.Cells.Item(81, 8).Value = 1320.2858
$ws.Cells.Item(81, 9).Value = 368.4
$ws.Cells.Item(81, 10).Value = 3700
$ws.Cells.Item(81, 11).Value = 736.8
$ws.Cells.Item(81, 12).Value = 7400
$ws.Cells.Item(81, 13).Value = 324.2
$ws.Cells.Item(81, 14).Value = -9522
$ws.Cells.Item(84, 8).Value = 1320.2858
$ws.Cells.Item(84, 9).Value = 368.4
$ws.Cells.Item(84, 10).Value = 3700
$ws.Cells.Item(84, 11).Value = 3684
$ws.Cells.Item(84, 12).Value = 37000
$ws.Cells.Item(84, 13).Value = 1620
$ws.Cells.Item(84, 14).Value = -47608
